# "support fast define enum" - adds a new demo sheet (Sheet2) that shows the
# fast-define-enum syntax ("myenum{A 0,B 14,C 23}:nameenum") plus a mirrored
# enum column (P) on the existing Sheet1 sample sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheet1: add column P ("myenum{A,B,C}:nameenum") to the first data block.
#    Cell write order below is deliberately chosen to match the original
#    authoring order (and therefore the shared-string table order): the data
#    value "A" was entered before the column header text was typed.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

$ws1.Cells.Item(3,16).Value = "A"
$ws1.Cells.Item(2,16).Value = "myenum{A,B,C}:nameenum"
$ws1.Cells.Item(4,16).Value = "A"
$ws1.Cells.Item(5,16).Value = "A"
$ws1.Cells.Item(6,16).Value = "A"
$ws1.Cells.Item(7,16).Value = "A"

# ---------------------------------------------------------------------------
# 2. Add the new "Sheet2" worksheet after Sheet1 and fill in the fast-enum
#    demo table.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2 = $wb.Worksheets.Add($null, $lastSheet)
$ws2.Name = "Sheet2"

# Row/column text is entered in an order that reproduces the original
# shared-string table order (new strings get interned in first-use order):
# "my:name", "D", "E", "F", "D 4", then finally the long fast-enum header.
$ws2.Cells.Item(1,1).Value = "st_level"
$ws2.Cells.Item(1,2).Value = "st_levelTable"
$ws2.Cells.Item(1,3).Value = "{""IsStringId"":false,""IsGenItemClass"":true,""JSONName"":""st_levelJSON"",""IsGenEnum"":true,""Path"":""toanstt/Resources/toandata"",""IsSeparatedJSON"":true}"

$ws2.Cells.Item(2,1).Value = "id"
$ws2.Cells.Item(2,3).Value = "my:name"

$ws2.Cells.Item(3,1).Value = 0
$ws2.Cells.Item(3,2).Value = "A"

$ws2.Cells.Item(6,3).Value = "D"

$ws2.Cells.Item(4,1).Value = 1
$ws2.Cells.Item(4,2).Value = "A"
$ws2.Cells.Item(4,3).Value = "E"

$ws2.Cells.Item(5,1).Value = 2
$ws2.Cells.Item(5,2).Value = "A"
$ws2.Cells.Item(5,3).Value = "F"

$ws2.Cells.Item(3,3).Value = "D 4"

$ws2.Cells.Item(2,2).Value = "myenum{A 0,B 14,C 23}:nameenum"

$ws2.Cells.Item(6,1).Value = 3
$ws2.Cells.Item(6,2).Value = "A"

$ws2.Cells.Item(7,1).Value = 4
$ws2.Cells.Item(7,2).Value = "A"
$ws2.Cells.Item(7,3).Value = "E"

# Column D on rows 3-7 carries the same (quote-prefixed/empty) style as
# column D on Sheet1 but no value.
for ($r = 3; $r -le 7; $r++) {
    $ws2.Cells.Item($r,4).Value = "'"
    $ws2.Cells.Item($r,4).Value = $null
}

$ws2.Columns.Item(2).ColumnWidth = 24.71

# ---------------------------------------------------------------------------
# 3. View state: Sheet2 becomes the active/selected tab with B2 highlighted;
#    Sheet1 keeps a plain A1:P7 selection anchored at G4.
# ---------------------------------------------------------------------------
$ws1.Select()
$excel.Application.Goto($ws1.Range("A1:P7"))
$ws1.Range("G4").Activate()

$ws2.Select()
$ws2.Range("B2").Select()
